$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 147, shifting existing rows 147:258 down to 148:259
$ws.Rows(147).Insert()

# Populate the newly inserted row 147 with the new record
$ws.Cells.Item(147, 1).Value  = 5
$ws.Cells.Item(147, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(147, 3).Value  = "Maule"
$ws.Cells.Item(147, 4).Value2 = 44907
$ws.Cells.Item(147, 5).Value  = 7
$ws.Cells.Item(147, 6).Value  = 100112024
$ws.Cells.Item(147, 7).Value  = "Choclo"
$ws.Cells.Item(147, 8).Value  = "Choclero"
$ws.Cells.Item(147, 9).Value  = "Primera"
$ws.Cells.Item(147, 10).Value = 30000
$ws.Cells.Item(147, 11).Value = 350
$ws.Cells.Item(147, 12).Value = 350
$ws.Cells.Item(147, 13).Value = 350
$ws.Cells.Item(147, 14).Value = "$/unidad"
$ws.Cells.Item(147, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(147, 16).Value = 350
$ws.Cells.Item(147, 17).Value = 1
$ws.Cells.Item(147, 18).Value = "Hortaliza"
